$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'sliding knee sleeve'
$ws.Cells.Item(2, 1).Value = 'volleyball spandex shorts'
$ws.Cells.Item(3, 1).Value = 'compression spandex men'
$ws.Cells.Item(4, 1).Value = 'volleyball knee pad'
$ws.Cells.Item(5, 1).Value = 'knee pads for volleyball'
$ws.Cells.Item(6, 1).Value = 'skateboard knee pads'
$ws.Cells.Item(7, 1).Value = 'starter youth compression pants'
$ws.Cells.Item(8, 1).Value = 'mens compression pants marvel'
$ws.Cells.Item(9, 1).Value = 'mouthguard basketball youth'
$ws.Cells.Item(10, 1).Value = 'nike compression leggings'
$ws.Cells.Item(11, 1).Value = 'knee pad sleeves'
$ws.Cells.Item(12, 1).Value = 'jordan flight mens basketball pants'
$ws.Cells.Item(13, 1).Value = 'athletic compression pants'
$ws.Cells.Item(14, 1).Value = 'nike pro compression leggings men'
$ws.Cells.Item(15, 1).Value = 'mens compression tights nike'
$ws.Cells.Item(16, 1).Value = 'ladies compression pants'
$ws.Cells.Item(17, 1).Value = 'mcdavid compression pants'
$ws.Cells.Item(18, 1).Value = 'basketball youth jersey'
$ws.Cells.Item(19, 1).Value = 'protec knee pads'
$ws.Cells.Item(20, 1).Value = 'skate knee pads'
$ws.Cells.Item(21, 1).Value = 'youth knee and elbow pads'
$ws.Cells.Item(22, 1).Value = 'under armour compression tights men'
$ws.Cells.Item(23, 1).Value = 'elbow and knee pads'
$ws.Cells.Item(24, 1).Value = 'adidas tights men'
$ws.Cells.Item(25, 1).Value = 'compression pants tesla'
$ws.Cells.Item(26, 1).Value = 'poc knee pads'
$ws.Cells.Item(27, 1).Value = 'basket knee pads'
$ws.Cells.Item(28, 1).Value = 'exercise knee pad'
$ws.Cells.Item(29, 1).Value = 'mens compression pants adidas'
$ws.Cells.Item(30, 1).Value = 'the rock mens basketball'
$ws.Cells.Item(31, 1).Value = 'men''s tights leggings'
$ws.Cells.Item(32, 1).Value = 'super compression leggings'
$ws.Cells.Item(33, 1).Value = 'under armour compression pants youth boys'
$ws.Cells.Item(34, 1).Value = 'knee pad compression pants'
$ws.Cells.Item(35, 1).Value = 'pantalones con rodilleras para hombre'
$ws.Cells.Item(36, 1).Value = 'basketball tights with knee pads for men'
$ws.Cells.Item(37, 1).Value = 'compression pants men basketball'
$ws.Cells.Item(38, 1).Value = 'knee leggings for basketball'
$ws.Cells.Item(39, 1).Value = 'pants with knee pads for men'
$ws.Cells.Item(40, 1).Value = 'knee pad pants for men'
$ws.Cells.Item(41, 1).Value = 'mens basketball compression pants'
$ws.Cells.Item(42, 1).Value = 'basketball yoga pants'
$ws.Cells.Item(43, 1).Value = 'basketball tights with pads'
$ws.Cells.Item(44, 1).Value = 'sliding pants mens'
$ws.Cells.Item(45, 1).Value = 'basketball compression pants knee pads'
$ws.Cells.Item(46, 1).Value = 'compression with knee pads'
$ws.Cells.Item(47, 1).Value = 'compression knee pads for men'
$ws.Cells.Item(48, 1).Value = 'men basketball compression knee pads'
$ws.Cells.Item(49, 1).Value = 'basketball leggings with pads'
$ws.Cells.Item(50, 1).Value = 'capri pads'
$ws.Cells.Item(51, 1).Value = 'basketball knee pad pants'
$ws.Cells.Item(52, 1).Value = 'tights with pads for men'
$ws.Cells.Item(53, 1).Value = 'mizuno slider knee pad'
$ws.Cells.Item(54, 1).Value = 'mizuno adult slider kneepad'
$ws.Cells.Item(55, 1).Value = 'compression knee leggings'
$ws.Cells.Item(56, 1).Value = '5 pad compression shorts'
$ws.Cells.Item(57, 1).Value = 'mizuno slider kneepad'
$ws.Cells.Item(58, 1).Value = 'youth basketball tights with knee pads'
$ws.Cells.Item(59, 1).Value = 'basketball knee pads for men'
$ws.Cells.Item(60, 1).Value = 'legging pads'
$ws.Cells.Item(61, 1).Value = 'compression capri men pack'
$ws.Cells.Item(62, 1).Value = 'softball knee sliding pad'
$ws.Cells.Item(63, 1).Value = 'basketball tights knee'
$ws.Cells.Item(64, 1).Value = 'compression pads basketball'
$ws.Cells.Item(65, 1).Value = 'leggings with knee pads'
$ws.Cells.Item(66, 1).Value = 'knee pads baseball'
$ws.Cells.Item(67, 1).Value = 'knee length tights men'
$ws.Cells.Item(68, 1).Value = 'basketball tight'
$ws.Cells.Item(69, 1).Value = 'baseball sliding knee pad'
$ws.Cells.Item(70, 1).Value = 'knee compression pants'
$ws.Cells.Item(71, 1).Value = 'baseball sliding pad'
$ws.Cells.Item(72, 1).Value = 'padded compression pants basketball'
$ws.Cells.Item(73, 1).Value = 'compression pants capri men'
$ws.Cells.Item(74, 1).Value = 'basketball men leggings'
$ws.Cells.Item(75, 1).Value = 'knee sliders softball'
$ws.Cells.Item(76, 1).Value = 'mens padded basketball tights'
$ws.Cells.Item(77, 1).Value = 'compression pants mens basketball'
$ws.Cells.Item(78, 1).Value = 'softball knee slider'
$ws.Cells.Item(79, 1).Value = 'mens weightlifting tights'
$ws.Cells.Item(80, 1).Value = 'mcdavid 6446 hex knee pads compression leg sleeve'
$ws.Cells.Item(81, 1).Value = 'softball sliding pad'
$ws.Cells.Item(82, 1).Value = 'pants pad'
$ws.Cells.Item(83, 1).Value = 'knee pads mens basketball'
$ws.Cells.Item(84, 1).Value = 'padded pants men'
$ws.Cells.Item(85, 1).Value = 'basketball legging'
$ws.Cells.Item(86, 1).Value = 'padded compression tights basketball'
$ws.Cells.Item(87, 1).Value = 'youth tights with knee pads'
$ws.Cells.Item(88, 1).Value = 'youth padded tights'
$ws.Cells.Item(89, 1).Value = 'baseball knee pads adult'
$ws.Cells.Item(90, 1).Value = 'mens capri leggings for sports'
$ws.Cells.Item(91, 1).Value = 'basketball tights with knee pads youth boys'
$ws.Cells.Item(92, 1).Value = 'youth compression pants with knee pads'
$ws.Cells.Item(93, 1).Value = 'mens capris pants'
$ws.Cells.Item(94, 1).Value = 'basketball knee pad leggings'
$ws.Cells.Item(95, 1).Value = 'padded pants for basketball'
$ws.Cells.Item(96, 1).Value = 'compression capri pants men'
$ws.Cells.Item(97, 1).Value = 'sliding knee pads'
$ws.Cells.Item(98, 1).Value = 'mens capri compression tights'
$ws.Cells.Item(99, 1).Value = 'compression below knee'
$ws.Cells.Item(100, 1).Value = 'softball sliding pads'
